$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing SEQUELIZE row (row 2) values ---
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = "Existen 5 (62%) elementos de JavaScript y 3 (37%) elementos de TypeScript"
$ws.Range("F2").Value = 8

# --- Add new BOOKSHELF row (row 3) ---
$ws.Range("A3").Value = "BOOKSHELF"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "Existen 2 (100%) elementos de JavaScript y 0 (0%) elementos de TypeScript"
$ws.Range("F3").Value = 2

# --- Add new PRISMA row (row 4) ---
$ws.Range("A4").Value = "PRISMA"
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = 7
$ws.Range("E4").Value = "Existen 0 (0%) elementos de JavaScript y 7 (100%) elementos de TypeScript"
$ws.Range("F4").Value = 7

# --- Apply thin box borders to every cell in rows 2-4 ---
$ws.Range("A2:F4").Borders.LineStyle = 1

# --- Dates (B3/B4) written as plain text "11/03/2023", matching B2 ---
# Mark the cells as Text first so Excel doesn't reinterpret the
# dd/mm/yyyy-looking string as a real date serial.
$ws.Range("B3:B4").NumberFormat = "@"
$ws.Range("B3").Value = "11/03/2023"
$ws.Range("B4").Value = "11/03/2023"
